$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Changed emergence check method - check last ten steps for emergence.
# As a result, the previously computed "A emerge percantage" (N) and
# "B emerge percantage" (O) values for every data row are reset to 0.
$ws.Range("N2:O52").Value = 0
